$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A forward by 3 days (re-run of the
# quarterly forecast model three days later).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 3
}

# Updated forecast values (column B) produced by the retrained model for
# rows 18 through 90 (the daylight hours where forecasted solar output
# changed).
$newB = @{
    18 = 17;   19 = 17;   20 = 19;   21 = 22;   22 = 108;
    23 = 116;  24 = 129;  25 = 148;  26 = 391;  27 = 417;
    28 = 448;  29 = 482;  30 = 898;  31 = 942;  32 = 988;
    33 = 1030; 34 = 1367; 35 = 1407; 36 = 1459; 37 = 1506;
    38 = 1792; 39 = 1819; 40 = 1855; 41 = 1880; 42 = 2087;
    43 = 2105; 44 = 2123; 45 = 2134; 46 = 2182; 47 = 2188;
    48 = 2193; 49 = 2198; 50 = 2168; 51 = 2170; 52 = 2168;
    53 = 2163; 54 = 2116; 55 = 2112; 56 = 2106; 57 = 2097;
    58 = 2003; 59 = 1991; 60 = 1972; 61 = 1951; 62 = 1779;
    63 = 1753; 64 = 1720; 65 = 1686; 66 = 1407; 67 = 1366;
    68 = 1317; 69 = 1276; 70 = 863;  71 = 815;  72 = 770;
    73 = 735;  74 = 348;  75 = 314;  76 = 289;  77 = 269;
    78 = 77;   79 = 62;   80 = 52;   81 = 45;   82 = 7;
    83 = 7;    84 = 6;    85 = 6;    86 = 5;    87 = 5;
    88 = 5;    89 = 5;    90 = 1
}

foreach ($r in $newB.Keys) {
    $ws.Cells.Item($r, 2).Value2 = $newB[$r]
}
